# Level 2 - Animation drop fix , Boss AI
#
# The "datetimeFigureOut" date field shown in the footer of the slide
# master and every slide layout was showing a stale cached date
# (28/03/2025). Bump it forward one day to 29/03/2025, matching what
# PowerPoint would cache the next time the deck is touched.

$p = $ppt.ActivePresentation
$newDate = "29/03/2025"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "28/03/2025") {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master footer date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's footer date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
